# Updates in June and July 2014
# - Adjust saved selections/scroll state on the first two sheets
# - Add a new "Biking-Maybe" sheet with a Tobler-hiking-function table
# - Bump the calc id so Excel knows to recalculate

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("Tobler002"): move the saved selection down to A46:B136 and
# scroll the view so row 116 is at the top.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("A46:B136").Select()
$excel.ActiveWindow.ScrollRow = 116
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------------
# Sheet 2 ("Sheet1"): move the saved selection to a single cell, F22.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("F22").Select()

# ---------------------------------------------------------------------------
# New sheet "Biking-Maybe": a 91-row, 2-column table.
#   Column A: the integers -45..45 (the "angle" input)
#   Column B, rows 1-46 (A = -45..0): formula  =1/COS((A*1.4)/180*PI())
#   Column B, rows 47-91 (A = 1..45): pasted/static decay values
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Biking-Maybe"

# Column A: -45 .. 45
$aArr = New-Object 'object[,]' 91,1
for ($i = 0; $i -lt 91; $i++) {
    $aArr[$i,0] = -45 + $i
}
$ws3.Range("A1:A91").Value = $aArr

# Column B, rows 1-46: formula =1/COS((A<row>*1.4)/180*PI())
for ($row = 1; $row -le 46; $row++) {
    $ws3.Cells.Item($row, 2).Formula = "=1/COS((A$row*1.4)/180*PI())"
}

# Column B, rows 47-91: static values (pasted-as-values data)
$staticVals = @(
    0.94073600387004974,
    0.88495127390790829,
    0.83241245027661315,
    0.78290506010728989,
    0.73623175718467948,
    0.69221074603919241,
    0.65067436899001063,
    0.61146783741132915,
    0.57444809084553139,
    0.53948276961858599,
    0.50644928837105219,
    0.47523399944234646,
    0.44573143636977447,
    0.41784362891561178,
    0.39147948203925742,
    0.3665542121076944,
    0.34298883440357325,
    0.32070969666113225,
    0.29964805394860622,
    0.27973968073271621,
    0.26092451641563863,
    0.24314634103564536,
    0.22635247817633528,
    0.21049352244206965,
    0.19552308913412081,
    0.18139758400764319,
    0.16807599120784555,
    0.15551967767810956,
    0.14369221250627179,
    0.13255919983050554,
    0.12208812406553485,
    0.11224820633533836,
    0.10301027111188625,
    0.09434662216243724,
    0.08623092700196516,
    0.07863810913372063,
    0.07154424744096648,
    0.06492648216765631,
    0.05876292699628209,
    0.05303258679822,
    0.047715280696534246,
    0.042791570144141346,
    0.038242691782221956,
    0.034050494905427824,
    0.030197383422318508
)
$bArr = New-Object 'object[,]' $staticVals.Count,1
for ($i = 0; $i -lt $staticVals.Count; $i++) {
    $bArr[$i,0] = $staticVals[$i]
}
$ws3.Range("B47:B91").Value = $bArr

# Select A1:B91 with the active cell at A91 (matches the saved view state),
# and leave this as the active sheet/tab so it is the one shown on open.
$ws3.Range("A1:B91").Select()
$ws3.Range("A91").Activate()
$ws3.Activate()

# ---------------------------------------------------------------------------
# Force a calc-id bump so downstream Excel knows a recalculation is needed.
# ---------------------------------------------------------------------------
$excel.CalculateFull()
